# Update "想去人数" (want-to-go count) figures that were refreshed by the
# gh-pages data regeneration (commit 456a3b4). The same underlying rows are
# duplicated across the "展览" sheet (1st sheet) and the "全部类型" sheet
# (4th sheet), so both need to be updated in lockstep.

$wb = $excel.ActiveWorkbook

# --- Sheet 1: 展览 ---------------------------------------------------------
$wsExpo = $wb.Worksheets.Item(1)
$wsExpo.Range("F2").Value  = 12631
$wsExpo.Range("F6").Value  = 276
$wsExpo.Range("F9").Value  = 12604
$wsExpo.Range("F10").Value = 17
$wsExpo.Range("F11").Value = 3121
$wsExpo.Range("F19").Value = 655
$wsExpo.Range("F21").Value = 6104
$wsExpo.Range("F23").Value = 3610

# --- Sheet 4: 全部类型 ------------------------------------------------------
$wsAll = $wb.Worksheets.Item(4)
$wsAll.Range("F2").Value  = 12631
$wsAll.Range("F6").Value  = 276
$wsAll.Range("F10").Value = 12604
$wsAll.Range("F11").Value = 17
$wsAll.Range("F12").Value = 3121
$wsAll.Range("F20").Value = 655
$wsAll.Range("F23").Value = 6104
$wsAll.Range("F25").Value = 3610
